$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the 6 new localization rows (key/value pairs) ---
# The sheet is a Key/Value table (col A = key, col B = value). We insert new
# rows at the appropriate spots, shifting everything below down, exactly the
# way Excel's Rows.Insert does it, then populate the new cells.

# 1) Two new rows before the existing "levelMatchTitle" row (orig row 13):
#    reveal / REVEAL, back / BACK
$ws.Rows("13:14").Insert()
$ws.Range("A13").Value = "reveal"
$ws.Range("B13").Value = "REVEAL"
$ws.Range("A14").Value = "back"
$ws.Range("B14").Value = "BACK"

# 2) Two new rows right after "levelMatchDesc" (now at row 16):
#    levelMatchNotFound / Climate does not match, try another location.
#    levelMatchFound / Climate Match Found!
$ws.Rows("17:18").Insert()
$ws.Range("A17").Value = "levelMatchNotFound"
$ws.Range("B17").Value = "Climate does not match, try another location."
$ws.Range("A18").Value = "levelMatchFound"
$ws.Range("B18").Value = "Climate Match Found!"

# 3) Two new rows right after "climateZone" (now at row 24):
#    zone / Zone, type / Type
$ws.Rows("25:26").Insert()
$ws.Range("A25").Value = "zone"
$ws.Range("B25").Value = "Zone"
$ws.Range("A26").Value = "type"
$ws.Range("B26").Value = "Type"

# Match the author's final selection state (last touched cell).
[void]$ws.Range("B17").Select()
